$wb = $excel.ActiveWorkbook

# --- Update "Ready for handoff" -> "Handed back: in sync with en-US" everywhere it appears ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = "Handed back: in sync with en-US"
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"

# --- Update the latest handback datetimes ---
$wsZh.Range("K2").Value = "2016-09-05 12:05:37"
$wsDe.Range("K2").Value = "2016-09-05 12:05:56"

# --- Clear the stale "handback file is not latest" error detail text ---
$wsZh.Range("P3").Value = ""
$wsDe.Range("P3").Value = ""

# --- Resize columns to match the report's refreshed layout ---
$wsOverview.Columns("E").ColumnWidth = 29.9777050018311
$wsOverview.Columns("F").ColumnWidth = 29.9777050018311

$wsZh.Columns("C").ColumnWidth = 29.9777050018311
$wsZh.Columns("P").ColumnWidth = 13.7470531463623

$wsDe.Columns("C").ColumnWidth = 29.9777050018311
$wsDe.Columns("P").ColumnWidth = 13.7470531463623

Write-Host "Report regenerated for handback"
